$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the three additional summary columns
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the header formatting already used by the rest of row 1 (bold,
# centered, thin border) by copying it from the last existing header cell.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values, row by row
$ws.Range("L2").Value = 90.1471541030692
$ws.Range("M2").Value = 211660
$ws.Range("N2").Value = 313.1065088757396

$ws.Range("L3").Value = 96.01937043276359
$ws.Range("M3").Value = 51893
$ws.Range("N3").Value = 336.9675324675325

$ws.Range("L4").Value = 88.30982863725519
$ws.Range("M4").Value = 147585
$ws.Range("N4").Value = 143.0087209302326

$ws.Range("L5").Value = 95.05584770392593
$ws.Range("M5").Value = 56061
$ws.Range("N5").Value = 159.7179487179487

$ws.Range("L6").Value = 17.70131305034959
$ws.Range("M6").Value = 1877
$ws.Range("N6").Value = 14.6640625

$ws.Range("L7").Value = 30.09544761111365
$ws.Range("M7").Value = 331
$ws.Range("N7").Value = 13.79166666666667
